$d = $word.ActiveDocument
$count = $d.Paragraphs.Count
$anchorPara = $d.Paragraphs.Item($count - 1)
$anchorRange = $anchorPara.Range
$anchorRange.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($count)
$insertRange = $newPara.Range

$xml = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:i/>
          <w:iCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:i/>
          <w:iCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>Backend side:</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:rStyle w:val="Hyperlink"/>
          <w:b/>
          <w:bCs/>
          <w:color w:val="auto"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:u w:val="none"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">Clone from git: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:i/>
          <w:iCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">git clone </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:i/>
          <w:iCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>git@github.com:orontamir</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:i/>
          <w:iCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>/</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:i/>
          <w:iCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>keshet</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:i/>
          <w:iCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>-server</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
          <w:color w:val="8EAADB" w:themeColor="accent1" w:themeTint="99"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">Insert into: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:i/>
          <w:iCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>Keshet</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:i/>
          <w:iCs/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>-server</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">Run the command: </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:i/>
          <w:iCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>npm</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:i/>
          <w:iCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> install --legacy-peer-deps</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:i/>
          <w:iCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:i/>
          <w:iCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>npm</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:i/>
          <w:iCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> run start</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
          <w:color w:val="8EAADB" w:themeColor="accent1" w:themeTint="99"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
          <w:color w:val="8EAADB" w:themeColor="accent1" w:themeTint="99"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>http://localhost:3000/search?q=</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="Arial" w:hint="cs"/>
          <w:i/>
          <w:iCs/>
          <w:color w:val="8EAADB" w:themeColor="accent1" w:themeTint="99"/>
          <w:rtl/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>אורון</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:i/>
          <w:iCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
    </w:p>
'@

[void]$insertRange.InsertXML($xml)

Write-Host ("Paragraph count after insert: " + $d.Paragraphs.Count)
